$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update existing Property/Value rows ---

# Version: 0.1.6 -> 0.1.7
$ws.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: updated publish date
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact (row 10): now shows the publisher-style contact line
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Contact (row 11): second contact, now a named person
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Make room for a new "Jurisdiction" row after the two Contact rows ---
# Shift rows 12-15 down to 13-16, carrying both formatting and values, working
# from the bottom up so nothing gets clobbered before it is copied.
function Copy-RowDown($srcRow, $dstRow) {
    $src = $ws.Range("A" + $srcRow + ":B" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":B" + $dstRow)
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $src.Copy()
    $dst.PasteSpecial(-4163)
}

Copy-RowDown 15 16
Copy-RowDown 14 15
Copy-RowDown 13 14
Copy-RowDown 12 13
$excel.CutCopyMode = 0

# Row 12 now holds the new "Jurisdiction" property (value left empty)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").ClearContents()
